$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-81)
# from serial date 45224 (2023-10-25) to 45233 (2023-11-03)
for ($r = 2; $r -le 81; $r++) {
    $ws.Cells.Item($r, 3).Value = 45233
}
